# ---------------------------------------------------------------------------
# Edit: "added in rest of data and added in WJ correlation testing"
#
# The underlying age table (Sheet1, columns A=ID, B=Age, C=Range) is extended:
#   - 8 new "6 to 8" rows are inserted in the middle of the table
#     (bkp059..bkp066), pushing what used to be the "Adult" placeholder
#     rows (bkp101..bkp118, previously Age=1) down.
#   - The "Adult" block grows from 18 rows (bkp101..bkp118) to 20 rows
#     (bkp101..bkp120), and every Adult row's Age is corrected from the
#     placeholder value 1 to the real value 18.
#   - The 8 new "6 to 8" rows (B38:B45) get WJ correlation formatting:
#     wrap text + an explicit black font color.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (row, ID, Age, Range) data for every row from 38 onward (rows 2-37
# are untouched by this edit).
$data = @(
    @(38, "bkp059", 6, "6 to 8"),
    @(39, "bkp060", 7, "6 to 8"),
    @(40, "bkp061", 6, "6 to 8"),
    @(41, "bkp062", 8, "6 to 8"),
    @(42, "bkp063", 8, "6 to 8"),
    @(43, "bkp064", 8, "6 to 8"),
    @(44, "bkp065", 7, "6 to 8"),
    @(45, "bkp066", 7, "6 to 8"),
    @(46, "bkp101", 18, "Adult"),
    @(47, "bkp102", 18, "Adult"),
    @(48, "bkp103", 18, "Adult"),
    @(49, "bkp104", 18, "Adult"),
    @(50, "bkp105", 18, "Adult"),
    @(51, "bkp106", 18, "Adult"),
    @(52, "bkp107", 18, "Adult"),
    @(53, "bkp108", 18, "Adult"),
    @(54, "bkp109", 18, "Adult"),
    @(55, "bkp110", 18, "Adult"),
    @(56, "bkp111", 18, "Adult"),
    @(57, "bkp112", 18, "Adult"),
    @(58, "bkp113", 18, "Adult"),
    @(59, "bkp114", 18, "Adult"),
    @(60, "bkp115", 18, "Adult"),
    @(61, "bkp116", 18, "Adult"),
    @(62, "bkp117", 18, "Adult"),
    @(63, "bkp118", 18, "Adult"),
    @(64, "bkp119", 18, "Adult"),
    @(65, "bkp120", 18, "Adult")
)

# New rows 56-65 do not exist yet in the sheet - stamp them with the same
# cell formatting as the existing "Adult" rows (ID col = style 1 / Arial,
# Range col = style 2 / wrap text) before writing values into them.
$ws.Range("A46:C46").Copy() | Out-Null
$ws.Range("A56:C65").PasteSpecial(-4122) | Out-Null
for ($r = 56; $r -le 65; $r++) {
    $ws.Rows.Item($r).RowHeight = 17
}

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# B38:B45 (the 8 new "6 to 8" rows) get the new WJ-correlation formatting:
# wrap text + explicit black font color (order matters to avoid spurious
# intermediate style/font entries: WrapText, then Font.Color, then Value).
for ($r = 38; $r -le 45; $r++) {
    $c = $ws.Cells.Item($r, 2)
    $c.WrapText = $true
    $c.Font.Color = 0
}

# Put the selection/scroll position where the author left it.
$ws.Range("D61").Select() | Out-Null
